$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (the "Förändrad" column) holds a date serial value that was
# bumped by one day (45633 -> 45634) for every data row (rows 2-34).
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45634
}
